$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Iterations=0, Success=False)
$ws.Range("D2").Value = 0.9999999999995148
$ws.Range("E2").Value = 0.9999999999995148

# Row 3 (Iterations=1, Success=True)
$ws.Range("D3").Value = 0.0007850321257758342
$ws.Range("E3").Value = 0.0007850321257758342

# Row 4 (Iterations=2, Success=True)
$ws.Range("D4").Value = 0.0001597317164635983
$ws.Range("E4").Value = 0.0001597317164635983

# Row 5 (Iterations=3, Success=True) - written without scientific notation
$ws.Range("D5").Value = 0.0000000000000007520521803276738
$ws.Range("E5").Value = 0.0000000000000007520521803276738

# Row 6 (Iterations=4, Success=False)
$ws.Range("D6").Value = 0.8954685934402525
$ws.Range("E6").Value = 0.8954685934402525

# Row 7 (Iterations=5, Success=True) - only F7 changes
$ws.Range("F7").Value = 5.102246284484863
